$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: fill in the previously-blank row with real data -----------
# Copy the formatting of row 16 (date/hours/activity columns) down into row 17
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A17").Value = 42513
$ws.Range("B17").Value = 2.5
$ws.Range("C17").Value = 'Drag&Drop angepasst. Versuch eine "SteineVorschau" einzubauen'

# --- Row 18: new data row, re-using the "highlighted" look of row 15 ---
$ws.Range("A15:C15").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A18").Value = 42514
$ws.Range("B18").Value = 4
$ws.Range("C18").Value = "Drag&Drop angepasst. Preview"

# --- Row 19: new blank spacer row (no border, default format) ----------
$ws.Range("A19").NumberFormat = "d-mmm"

$ws.Range("B19").HorizontalAlignment = -4108   # xlCenter (touch ...
$ws.Range("B19").HorizontalAlignment = 1       # ... then back to xlGeneral)
$ws.Range("C19").HorizontalAlignment = 1

# --- Row 20: the running total, now covering the two new rows ----------
$ws.Range("B20").NumberFormat = "0.0"
$ws.Range("B20").Formula = "=SUM(B3:B18)"

# --- Sheet bookkeeping ---------------------------------------------------
$ws.Range("C11").Select()

Write-Host "done"
